$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 39: Date, Effort[h], (no Additional Effort), Comment
$ws.Range("A39").NumberFormat = "ddd\ dd/mm/yyyy"
$ws.Range("A39").Value = 41221

$ws.Range("B39").Value = 3

$ws.Range("D39").Value = "Installer creation scripts continued, missing readMe files added"

# Update selection to match the new active cell
$ws.Range("B39").Select()
